$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D column (and let E recalc via existing formulas) for rows 2-6
$ws.Range("D2").Value = 261656214
$ws.Range("D3").Value = 174187509
$ws.Range("D4").Value = 610382383
$ws.Range("D5").Value = 164000687
$ws.Range("D6").Value = 549534067

# Rows 7 and 8: team names swap (Zakka_S2l <-> Mkm_s2l) along with their C/D values
$ws.Range("B7").Value = "Mkm_s2l"
$ws.Range("C7").Value = 3202314817
$ws.Range("D7").Value = 278322637

$ws.Range("B8").Value = "Zakka_S2l"
$ws.Range("C8").Value = 5934351328
$ws.Range("D8").Value = 514764915

# Update the selection to reflect the new active cell / selected range
$ws.Range("E2:E8").Select() | Out-Null
